$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking columns (D = Price, E = Volume(1h)) must be forced to
# Text so Excel does not auto-coerce them into Number/Percentage values,
# matching the original inline-string cell type. ClearFormats() afterwards
# restores the default (unstyled) cell formatting.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}


# Row 2
Set-TextValue $ws.Range("D2") "308.75"
Set-TextValue $ws.Range("E2") "1.18%"

# Row 3
Set-TextValue $ws.Range("D3") "38.48"
Set-TextValue $ws.Range("E3") "7.60%"

# Row 4
Set-TextValue $ws.Range("D4") "5.103"
Set-TextValue $ws.Range("E4") "1.50%"

# Row 5
Set-TextValue $ws.Range("D5") "0.08122"
Set-TextValue $ws.Range("E5") "1.03%"

# Row 6
Set-TextValue $ws.Range("D6") "1.970"
Set-TextValue $ws.Range("E6") "5.04%"

# Row 7
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D7") "4.205"
Set-TextValue $ws.Range("E7") "1.42%"

# Row 8
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws.Range("D8") "7.933"
Set-TextValue $ws.Range("E8") "1.91%"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D9") "0.9292"
Set-TextValue $ws.Range("E9") "0.89%"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D10") "0.1449"
Set-TextValue $ws.Range("E10") "13.03%"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1962"
Set-TextValue $ws.Range("E11") "2.53%"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.09102"
Set-TextValue $ws.Range("E12") "0.31%"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.03508"
Set-TextValue $ws.Range("E13") "1.16%"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.09811"
Set-TextValue $ws.Range("E14") "-0.51%"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001407"
Set-TextValue $ws.Range("E15") "-0.56%"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D16") "0.006100"
Set-TextValue $ws.Range("E16") "-2.79%"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.655"
Set-TextValue $ws.Range("E17") "-4.24%"

# Row 18
Set-TextValue $ws.Range("D18") "3.487"
Set-TextValue $ws.Range("E18") "2.71%"

# Row 19
Set-TextValue $ws.Range("D19") "0.3445"
Set-TextValue $ws.Range("E19") "0.78%"

# Row 20
Set-TextValue $ws.Range("D20") "0.1333"
Set-TextValue $ws.Range("E20") "0.89%"

# Row 21
Set-TextValue $ws.Range("D21") "4.802"
Set-TextValue $ws.Range("E21") "-7.31%"

# Row 22
Set-TextValue $ws.Range("E22") "6.34%"

# Row 23
Set-TextValue $ws.Range("D23") "0.04433"
Set-TextValue $ws.Range("E23") "0.43%"

# Row 24
Set-TextValue $ws.Range("D24") "0.001219"
Set-TextValue $ws.Range("E24") "-1.28%"

# Row 25
Set-TextValue $ws.Range("D25") "0.004835"
Set-TextValue $ws.Range("E25") "4.85%"

# Row 27
Set-TextValue $ws.Range("D27") "0.0001302"
Set-TextValue $ws.Range("E27") "3.97%"

# Row 39
Set-TextValue $ws.Range("E39") "8.46%"

# Row 40
Set-TextValue $ws.Range("D40") "0.05145"
Set-TextValue $ws.Range("E40") "-2.47%"

# Row 41
Set-TextValue $ws.Range("D41") "0.007462"
Set-TextValue $ws.Range("E41") "-2.01%"

# Row 42
Set-TextValue $ws.Range("D42") "0.01013"
Set-TextValue $ws.Range("E42") "0.02%"

# Row 43
Set-TextValue $ws.Range("D43") "0.1363"
Set-TextValue $ws.Range("E43") "0.87%"

# Row 44
Set-TextValue $ws.Range("D44") "0.002143"
Set-TextValue $ws.Range("E44") "-0.49%"

# Row 45
Set-TextValue $ws.Range("D45") "0.009211"
Set-TextValue $ws.Range("E45") "-4.16%"

# Row 46
Set-TextValue $ws.Range("D46") "0.00006278"
Set-TextValue $ws.Range("E46") "2.62%"

# Row 47
Set-TextValue $ws.Range("E47") "-0.03%"

# Row 48
Set-TextValue $ws.Range("D48") "0.003057"

# Row 49
Set-TextValue $ws.Range("D49") "0.001600"
Set-TextValue $ws.Range("E49") "-3.63%"

# Row 50
Set-TextValue $ws.Range("D50") "0.00002103"
Set-TextValue $ws.Range("E50") "-0.03%"

# Row 51
Set-TextValue $ws.Range("D51") "0.0002003"
Set-TextValue $ws.Range("E51") "-0.03%"
